# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps for the handed-off rows, and marks their "Priority" column as
# "ht" (handoff type) on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 11, 12, 13, 14)

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-24 16:21:24"
}

# --- zh-cn sheet: column H = "Latest Handoff Datetime", column E = "Priority" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-08-24 16:21:19"
    $wsZhCn.Range("E$r").Value = "ht"
}

# --- de-de sheet: column H = "Latest Handoff Datetime", column E = "Priority" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("H$r").Value = "2016-08-24 16:21:24"
    $wsDeDe.Range("E$r").Value = "ht"
}
